$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (changed) date column C was updated from 2023-10-05
# (serial 45204) to 2023-10-06 (serial 45205) for every data row (2-29).
$ws.Range("C2:C29").Value = 45205
